# daily auto push: 2025-10-07 01:59 UTC
# Append a new daily-ranking row (row 74) to the sheet, mirroring the
# existing layout: 日付 (date, text), 曜日 (weekday, text), 時刻 (hour, number),
# ランキング (ranking, number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the previous row (r73: date "2025/10/07", weekday "火", hour 6,
# ranking 201) down into the new row 74. This keeps the date/weekday cells
# stored as plain text (matching the rest of the sheet) instead of letting
# Excel auto-convert the date-like string into a real date value/format.
$ws.Range("A73:D73").Copy($ws.Range("A74:D74"))

# The new entry shares the same date (2025/10/07) and weekday (火) as row
# 73, but was logged at hour 10; the ranking stays 201.
$ws.Cells.Item(74, 3).Value = 10
